# The presentation's live theme (ppt/theme/theme2.xml, the "Integral" /
# "Red Violet" theme used by the slide master) is swapped back to the
# original Office-default colour scheme that previously lived in
# ppt/theme/theme1.xml ("Office Theme" / "Office").
#
# PowerPoint's automation model exposes the twelve theme colour scheme
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) through
# Slide.ThemeColorScheme, in that fixed order. Writing to them edits the
# underlying theme part's <a:clrScheme> colour values directly, without
# touching anything else in the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$tcs.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
